$d = $word.ActiveDocument

# --- Paragraph 1: title text change ---
$d.Content.Find.Execute("Test document", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Output of test.htm", 2)

# --- Paragraph 2: replace with the "rendered" contents of test.htm, one
#     line of markup per segment, separated by manual line breaks (so the
#     OOXML ends up as a single run with alternating <w:t>/<w:br/>). ---
$manualBreak = [char]11
$lines = @(
    "<!DOCTYPE html>",
    "<html>",
    "<head>",
    "<title>test page</title>",
    "</head>",
    "<body>",
    "",
    "<h1>This is a test page</h1>",
    "<p>This is a test page, that will allow me to test my html to word doc convert</p>",
    "</body>",
    "</html>",
    "",
    ""
)
$htmlText = $lines -join $manualBreak

$p2 = $d.Paragraphs(2).Range
$p2.Text = $htmlText

Write-Output "done"
